# Insert a new data row above row 403 (shifting existing rows 403:456 down
# to 404:457) and populate the new row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(403).Insert()

$ws.Range("A403").Value = 6
$ws.Range("B403").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C403").Value = "Metropolitana"
$ws.Range("D403").Value = 44748
$ws.Range("E403").Value = 13
$ws.Range("F403").Value = 100112043
$ws.Range("G403").Value = "Pepino ensalada"
$ws.Range("H403").Value = "Sin especificar"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 300
$ws.Range("K403").Value = 17000
$ws.Range("L403").Value = 19000
$ws.Range("M403").Value = 18200
$ws.Range("N403").Value = "`$/caja 60 unidades"
$ws.Range("O403").Value = "Región de Arica y Parinacota"
$ws.Range("P403").Value = 303
$ws.Range("Q403").Value = 60
$ws.Range("R403").Value = "Hortaliza"
